$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume/Change (E) columns
# Column D values are stored as text (to preserve formats like "30.363.24"),
# so force Text number format before assigning to avoid Excel auto-converting
# numeric-looking strings into actual numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.363.24'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.866.79'
$ws.Range("E3").Value = '  +0.27%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.60'
$ws.Range("E5").Value = '  -0.96%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4705'
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("E8").Value = '  -1.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06567'
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.25'
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07879'
$ws.Range("E11").Value = '  -0.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '96.96'
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.875.19'
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6950'
$ws.Range("E14").Value = '  +1.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.099'
$ws.Range("E15").Value = '  -1.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '268.02'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.407.11'
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007636'
$ws.Range("E19").Value = '  +3.57%  '
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.126.68'
$ws.Range("E21").Value = '  +0.57%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.226'
$ws.Range("E23").Value = '  -1.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.174'
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.406'
$ws.Range("E25").Value = '  +2.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.34'
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.90'
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.947'
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.363'
$ws.Range("E29").Value = '  -2.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09922'
$ws.Range("E30").Value = '  +1.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.370'
$ws.Range("E32").Value = '  -0.81%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.055'
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04747'
$ws.Range("E34").Value = '  +1.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.133'
$ws.Range("E35").Value = '  +0.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7013'
$ws.Range("E36").Value = '  -0.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.715'
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("E38").Value = '  +0.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.793'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.301'
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.04'
$ws.Range("E41").Value = '  -1.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.951'
$ws.Range("E42").Value = '  +1.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8368'
$ws.Range("E45").Value = '  -0.99%  '
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '966.46'
$ws.Range("E47").Value = '  +1.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.125'
$ws.Range("E48").Value = '  -0.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.176'
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.49'
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05679'
$ws.Range("E51").Value = '  +0.43%  '
